$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet view: zoom + selection ---
$win = $excel.ActiveWindow
$win.Zoom = 40
$ws.Range("L9").Select() | Out-Null

# --- Header text for the third block (K1:N1) ---
$ws.Range("M1").Value2 = "exercices"
$ws.Range("N1").Value2 = "max exercices"

# --- Row 2 ---
$ws.Range("M2").Value2 = 11
$ws.Range("N2").Value2 = 700

# --- Row 3 ---
$ws.Range("C3").Value2 = 41
$ws.Range("M3").Value2 = 14
$ws.Range("N3").Value2 = 500

# --- Row 4 ---
$ws.Range("N4").Value2 = 650

# --- Row 5 ---
$ws.Range("N5").Value2 = 450

# --- Row 6 ---
$ws.Range("N6").Value2 = 450

# --- Row 7 ---
$ws.Range("M7").Value2 = 1
$ws.Range("N7").Value2 = 675

# --- Row 8 ---
$ws.Range("N8").Value2 = 680
